$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet stores every data cell as plain text (t="inlineStr" in the
# original file), including numeric-looking prices/counters and percentage
# strings such as "5.74%". Force each touched cell to a text number format
# first so Excel does not auto-convert the new values into numbers or
# percentages, which would change the cells underlying type.
$updates = @(
    @{ Cell = "D2"; Value = "304.64" }
    @{ Cell = "E2"; Value = "5.74%" }
    @{ Cell = "G2"; Value = "7" }
    @{ Cell = "D3"; Value = "34.89" }
    @{ Cell = "E3"; Value = "12.42%" }
    @{ Cell = "G3"; Value = "7" }
    @{ Cell = "D4"; Value = "5.145" }
    @{ Cell = "E4"; Value = "4.31%" }
    @{ Cell = "G4"; Value = "7" }
    @{ Cell = "D5"; Value = "0.07764" }
    @{ Cell = "E5"; Value = "5.46%" }
    @{ Cell = "G5"; Value = "7" }
    @{ Cell = "D6"; Value = "2.353" }
    @{ Cell = "E6"; Value = "4.12%" }
    @{ Cell = "G6"; Value = "7" }
    @{ Cell = "D7"; Value = "8.018" }
    @{ Cell = "E7"; Value = "3.69%" }
    @{ Cell = "G7"; Value = "7" }
    @{ Cell = "D8"; Value = "3.938" }
    @{ Cell = "E8"; Value = "5.31%" }
    @{ Cell = "G8"; Value = "7" }
    @{ Cell = "D9"; Value = "0.9228" }
    @{ Cell = "E9"; Value = "1.80%" }
    @{ Cell = "G9"; Value = "7" }
    @{ Cell = "D10"; Value = "0.1002" }
    @{ Cell = "E10"; Value = "15.31%" }
    @{ Cell = "G10"; Value = "7" }
    @{ Cell = "D11"; Value = "0.1799" }
    @{ Cell = "E11"; Value = "6.84%" }
    @{ Cell = "G11"; Value = "7" }
    @{ Cell = "D12"; Value = "0.08541" }
    @{ Cell = "E12"; Value = "3.76%" }
    @{ Cell = "G12"; Value = "7" }
    @{ Cell = "D13"; Value = "0.03304" }
    @{ Cell = "E13"; Value = "6.08%" }
    @{ Cell = "G13"; Value = "7" }
    @{ Cell = "D14"; Value = "0.09896" }
    @{ Cell = "E14"; Value = "-0.37%" }
    @{ Cell = "G14"; Value = "7" }
    @{ Cell = "D15"; Value = "0.001482" }
    @{ Cell = "E15"; Value = "-1.36%" }
    @{ Cell = "G15"; Value = "7" }
    @{ Cell = "D16"; Value = "0.005746" }
    @{ Cell = "E16"; Value = "-0.26%" }
    @{ Cell = "G16"; Value = "7" }
    @{ Cell = "D17"; Value = "3.469" }
    @{ Cell = "E17"; Value = "-0.61%" }
    @{ Cell = "G17"; Value = "7" }
    @{ Cell = "D18"; Value = "2.133" }
    @{ Cell = "E18"; Value = "1.73%" }
    @{ Cell = "G18"; Value = "7" }
    @{ Cell = "E19"; Value = "1.24%" }
    @{ Cell = "G19"; Value = "7" }
    @{ Cell = "D20"; Value = "0.1301" }
    @{ Cell = "E20"; Value = "0.57%" }
    @{ Cell = "G20"; Value = "7" }
    @{ Cell = "D21"; Value = "4.303" }
    @{ Cell = "E21"; Value = "12.45%" }
    @{ Cell = "G21"; Value = "7" }
    @{ Cell = "E22"; Value = "12.30%" }
    @{ Cell = "G22"; Value = "7" }
    @{ Cell = "D23"; Value = "0.04563" }
    @{ Cell = "E23"; Value = "0.13%" }
    @{ Cell = "G23"; Value = "7" }
    @{ Cell = "D24"; Value = "0.001217" }
    @{ Cell = "E24"; Value = "0.59%" }
    @{ Cell = "G24"; Value = "7" }
    @{ Cell = "D25"; Value = "0.004461" }
    @{ Cell = "E25"; Value = "7.75%" }
    @{ Cell = "G25"; Value = "7" }
    @{ Cell = "D26"; Value = "0.0001298" }
    @{ Cell = "E26"; Value = "-0.24%" }
    @{ Cell = "G26"; Value = "7" }
    @{ Cell = "E27"; Value = "8.76%" }
    @{ Cell = "G27"; Value = "7" }
    @{ Cell = "G28"; Value = "7" }
    @{ Cell = "G29"; Value = "7" }
    @{ Cell = "G30"; Value = "7" }
    @{ Cell = "G31"; Value = "7" }
    @{ Cell = "G32"; Value = "7" }
    @{ Cell = "G33"; Value = "7" }
    @{ Cell = "G34"; Value = "7" }
    @{ Cell = "G35"; Value = "7" }
    @{ Cell = "G36"; Value = "7" }
    @{ Cell = "G37"; Value = "7" }
    @{ Cell = "G38"; Value = "7" }
    @{ Cell = "D39"; Value = "0.01789" }
    @{ Cell = "E39"; Value = "13.51%" }
    @{ Cell = "G39"; Value = "7" }
    @{ Cell = "D40"; Value = "0.04760" }
    @{ Cell = "E40"; Value = "6.63%" }
    @{ Cell = "G40"; Value = "7" }
    @{ Cell = "D41"; Value = "0.007767" }
    @{ Cell = "E41"; Value = "5.68%" }
    @{ Cell = "G41"; Value = "7" }
    @{ Cell = "D42"; Value = "0.1413" }
    @{ Cell = "E42"; Value = "6.68%" }
    @{ Cell = "G42"; Value = "7" }
    @{ Cell = "D43"; Value = "0.007071" }
    @{ Cell = "E43"; Value = "-26.46%" }
    @{ Cell = "G43"; Value = "7" }
    @{ Cell = "D44"; Value = "0.002103" }
    @{ Cell = "E44"; Value = "-6.18%" }
    @{ Cell = "G44"; Value = "7" }
    @{ Cell = "D45"; Value = "0.009514" }
    @{ Cell = "E45"; Value = "12.75%" }
    @{ Cell = "G45"; Value = "7" }
    @{ Cell = "D46"; Value = "0.00006115" }
    @{ Cell = "E46"; Value = "0.00%" }
    @{ Cell = "G46"; Value = "7" }
    @{ Cell = "E47"; Value = "-0.16%" }
    @{ Cell = "G47"; Value = "7" }
    @{ Cell = "D48"; Value = "2.732" }
    @{ Cell = "E48"; Value = "24.89%" }
    @{ Cell = "G48"; Value = "7" }
    @{ Cell = "E49"; Value = "-0.22%" }
    @{ Cell = "G49"; Value = "7" }
    @{ Cell = "D50"; Value = "0.00002099" }
    @{ Cell = "E50"; Value = "-0.16%" }
    @{ Cell = "G50"; Value = "7" }
    @{ Cell = "E51"; Value = "-0.16%" }
    @{ Cell = "G51"; Value = "7" }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.NumberFormat = "@"
    $range.Value = $u.Value
}

Write-Host "Updated $($updates.Count) cells with refreshed symbol data"
